$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing backup codes (A2:A4) with new values
$ws.Range("A2").Value = "QS3W554CY3ZX"
$ws.Range("A3").Value = "Q8YK0WWBPD6H"
$ws.Range("A4").Value = "YYZGQ1P8K4EG"

# Clear old trailing codes that used to live at A14:A17
$ws.Range("A14:A17").ClearContents()

# Write the new set of codes into rows 11-16
$ws.Range("A11").Value = "3Z6ADAYX8TXX"
$ws.Range("A12").Value = "96R6XXZ5H6HD"
$ws.Range("A13").Value = "51XZBTACEQGM"
$ws.Range("A14").Value = "T91KP1V5CTSN"
$ws.Range("A15").Value = "YW0AXXBBWGJP"
$ws.Range("A16").Value = "KBN7HS57G4H3"

# Restore the active cell selection to A4
$ws.Range("A4").Select()
